$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 and Row 4 have their data (Fecha, Calidad, Volumen, Precio mínimo,
# Precio máximo, Precio promedio ponderado, Precio $/Kg) swapped.

# New values for row 2 (previously row 4's values)
$ws.Range("D2").Value = 44923
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 80
$ws.Range("N2").Value = 7500
$ws.Range("O2").Value = 8000
$ws.Range("P2").Value = 7625
$ws.Range("S2").Value = 7625

# New values for row 4 (previously row 2's values)
$ws.Range("D4").Value = 44881
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 11250
$ws.Range("O4").Value = 11250
$ws.Range("P4").Value = 11250
$ws.Range("S4").Value = 11250
